$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: A=1 (unchanged), B and C updated
$ws.Range("B3").Value = 2436
$ws.Range("C3").Value = 50

# Row 4: A=2 (unchanged), B and C updated
$ws.Range("B4").Value = 3033
$ws.Range("C4").Value = 180

# Row 5: A=3 (unchanged), B and C updated
$ws.Range("B5").Value = 3035
$ws.Range("C5").Value = 32

# Row 6: A=4 (unchanged), B updated, C stays 27
$ws.Range("B6").Value = 3141
$ws.Range("C6").Value = 27

# Row 7: A=5 (unchanged), B and C updated
$ws.Range("B7").Value = 3189
$ws.Range("C7").Value = 27

# Row 8: A changes from 8 to 7, B and C get new values
$ws.Range("A8").Value = 7
$ws.Range("B8").Value = 3588
$ws.Range("C8").Value = 35

# Row 9: brand new row; copy formatting from A8 (bold/border/centered style)
$ws.Range("A8").Copy()
$ws.Range("A9").PasteSpecial(-4122)
$ws.Range("A9").Value = 9
$ws.Range("B9").Value = 6104
$ws.Range("C9").Value = 36

# Row 10: brand new row; A=11, B=6411 (old row8 B), C=26 (old row8 C)
$ws.Range("A8").Copy()
$ws.Range("A10").PasteSpecial(-4122)
$ws.Range("A10").Value = 11
$ws.Range("B10").Value = 6411
$ws.Range("C10").Value = 26

$excel.CutCopyMode = $false
